# Update gh-pages to output generated at 456a3b4
# Refreshes scraped interest-count ("想去人数") and min-price ("最低票价")
# figures across the 展览 / 演出 / 本地生活 / 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Values = @{
    "F4"  = 10350
    "F5"  = 767
    "F8"  = 489
    "F9"  = 450
    "F11" = 280
    "F13" = 12940
    "F14" = 12940
    "F18" = 202
    "F20" = 192
    "F21" = 2778
    "F23" = 106
    "F25" = 133
    "G27" = 80
    "F29" = 2180
    "F30" = 1158
    "F31" = 4348
    "F33" = 3934
    "F34" = 966
    "F35" = 2685
    "F38" = 1410
    "F39" = 221
    "F42" = 168
    "F43" = 620
    "F44" = 881
    "F47" = 365
    "F49" = 202
    "F50" = 235
}
foreach ($cellRef in $sheet1Values.Keys) {
    $ws1.Range($cellRef).Value = $sheet1Values[$cellRef]
}

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Values = @{
    "F5"  = 72
    "F10" = 82
    "F18" = 47
}
foreach ($cellRef in $sheet2Values.Keys) {
    $ws2.Range($cellRef).Value = $sheet2Values[$cellRef]
}

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("G2").Value = 88

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Values = @{
    "F4"  = 10350
    "F5"  = 767
    "F7"  = 489
    "F8"  = 450
    "F10" = 280
    "F11" = 12940
    "F12" = 12940
    "G13" = 88
    "F17" = 202
    "F19" = 2778
    "F22" = 133
    "G24" = 80
    "F26" = 2180
    "F27" = 1158
    "F31" = 4349
    "F32" = 3934
    "F33" = 966
    "F34" = 2685
    "F38" = 47
    "F39" = 221
    "F42" = 620
    "F44" = 881
    "F47" = 365
    "F49" = 202
    "F50" = 235
}
foreach ($cellRef in $sheet4Values.Keys) {
    $ws4.Range($cellRef).Value = $sheet4Values[$cellRef]
}

Write-Output "applied gh-pages data refresh (456a3b4)"
